# Apply the edits described by the commit "fixed trees section and added 2 more programs."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. "fixed trees section": rows 5 and 6 had the Workforce-Development-Week
#    and GED figures swapped. Correct the G/H pair on each row.
# ---------------------------------------------------------------------------
$ws.Range("G5").Value = "Workforce Development Week"
$ws.Range("H5").Value = 105

$ws.Range("G6").Value = "GED"
$ws.Range("H6").Value = 2

# ---------------------------------------------------------------------------
# 2. "added 2 more programs": a new Stat/Category table in columns S:T.
# ---------------------------------------------------------------------------
$ws.Range("S1").Value = "Stat"
$ws.Range("T1").Value = "Category"

$ws.Range("S2").Value = 540
$ws.Range("T2").Value = "families' lights on "

$ws.Range("S3").Value = 900000
$ws.Range("S3").NumberFormat = "#,##0"
$ws.Range("T3").Value = "rental and utility assistance"

$ws.Range("S4").Value = 55
$ws.Range("T4").Value = "didn’t not qualify for government support"

$ws.Range("S5").Value = 1290
$ws.Range("T5").Value = "client records "

$ws.Range("S6").Value = 3999
$ws.Range("T6").Value = "residents "

$ws.Range("S7").Value = 95
$ws.Range("T7").Value = "at 70% of below AMI"

$ws.Range("S8").Value = 79.8
$ws.Range("T8").Value = "Hispanic "

$ws.Range("S9").Value = 74.3
$ws.Range("T9").Value = "Spanish Speakers "

$ws.Range("S10").Value = 95
$ws.Range("T10").Value = "at 70% or below median income "

# Widen the new Category column so the longer text fits.
$ws.Columns.Item(20).ColumnWidth = 24.5

# Restore the active-cell selection recorded in the saved workbook.
[void]$ws.Range("E8").Select()
